$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last row (48) currently carries the "current day" number format
# (date-only, style 3). Reset it to the normal "date + time" format
# (style 2) used by all the other data rows, since row 49 will become
# the new last row.
$ws.Cells.Item(48, 1).NumberFormat = $ws.Cells.Item(47, 1).NumberFormat

# Append the new daily-update row (row 49).
$ws.Cells.Item(49, 1).Value = 45789
$ws.Cells.Item(49, 2).Value = 199
$ws.Cells.Item(49, 3).Value = 211
$ws.Cells.Item(49, 4).Value = 208

# Give the new last row the "current day" (date-only) number format that
# row 48 used to have.
$ws.Cells.Item(49, 1).NumberFormat = "YYYY-MM-DD"
